$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 180: admissions corrections
$ws.Range("G180").Value = 6
$ws.Range("I180").Value = 1

# Row 181: new cases / hospitalises hors SI corrections
$ws.Range("C181").Value = 15
$ws.Range("G181").Value = 5

# Row 182: new cases / hospitalises hors SI corrections
$ws.Range("C182").Value = 10
$ws.Range("G182").Value = 4
$ws.Range("I182").Value = 1

# Row 183: fill in the day's figures (previously blank)
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 1
$ws.Range("F183").Value = 1
$ws.Range("G183").Value = 3
$ws.Range("I183").Value = 1
$ws.Range("L183").Value = "0"
$ws.Range("M183").Value = "0"

# Update the last-clicked cell to match the author's session
$ws.Range("P173").Select()
